$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the SetBalance method result cached in the sheet: F6 balance 1000 -> 1008
$ws.Range("F6").Value = 1008

# Append three more user rows (90-92) by cloning the existing
# "moses / bro / 1234 / m@g.c / Male / 0" template row - the same shape
# already used for the run of rows 68-89 - so the new rows keep identical
# cell types/formatting (copy+paste preserves the text-typed "1234" in
# column C instead of letting it auto-convert to a number).
foreach ($r in 90..92) {
    $ws.Range("A68:F68").Copy()
    $ws.Range(("A{0}:F{0}" -f $r)).PasteSpecial()
}

# Scroll the view down to where the new rows were added, keeping the
# original F67 selection.
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("F67").Select()
